# Update IPSA worksheet with latest data (Actualizar 03-08-2021 13-53-03)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the tail of the data series ---
# Replace the stale duplicated rows (794:803) with the three new observations
# that continue the series (03/03/2021, 03/04/2021, 03/05/2021).
$ws.Range("A794").Value = 44258
$ws.Range("B794").Value = 4787.45

$ws.Range("A795").Value = 44259
$ws.Range("B795").Value = 4700.57

$ws.Range("A796").Value = 44260
$ws.Range("B796").Value = 4713.732

# The remaining previously-filled rows become blank placeholders again.
$ws.Range("A797:B803").ClearContents()

# Drop the extra trailing blank rows that are no longer part of the sheet.
$ws.Range("A805:A812").EntireRow.Delete()

# --- Cosmetic / layout updates that come along with the refresh ---
# Column B widens to fit the (unwrapped) header text.
$ws.Columns("B").ColumnWidth = 79

# Row 1 no longer needs the tall custom height now that wrapping isn't forced.
$ws.Rows("1").AutoFit()

# Keep the named range in sync with the new data extent.
$wb.Names.Item(1).RefersTo = "=IPSA!`$A`$1:`$B`$795"

# Restore the view/selection state to match where the user left off.
$win = $excel.ActiveWindow
$win.ScrollRow = 794
$win.ScrollColumn = 2
$ws.Range("B798").Select()
